# Apply updated TPM values for Csf3-Csf3r LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Resolving-Mac" sending-cluster block (rows 11-13); data now ends at row 10
$ws.Rows("11:13").Delete()

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5070163333333334
$ws.Range("H2").Value = 1.521049
$ws.Range("I2").Value = 0.697371256392364
$ws.Range("J2").Value = 0.697371256392364
$ws.Range("M2").Value = 0.07904133333333334
$ws.Range("N2").Value = 0.237124
$ws.Range("O2").Value = 0.0007343710751920149
$ws.Range("P2").Value = 0.0007343710751920148
$ws.Range("Q2").Value = 0.04007524700844445
$ws.Range("R2").Value = 0.360677223076
$ws.Range("S2").Value = 0.0005121292793648666
$ws.Range("T2").Value = 0.0005121292793648666

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5070163333333334
$ws.Range("H3").Value = 1.521049
$ws.Range("I3").Value = 0.697371256392364
$ws.Range("J3").Value = 0.697371256392364
$ws.Range("M3").Value = 68.78716633333333
$ws.Range("N3").Value = 206.361499
$ws.Range("O3").Value = 0.6390998629361258
$ws.Range("P3").Value = 0.6390998629361259
$ws.Range("Q3").Value = 34.87621685471678
$ws.Range("R3").Value = 313.885951692451
$ws.Range("S3").Value = 0.4456898743759536
$ws.Range("T3").Value = 0.4456898743759537

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5070163333333334
$ws.Range("H4").Value = 1.521049
$ws.Range("I4").Value = 0.697371256392364
$ws.Range("J4").Value = 0.697371256392364
$ws.Range("M4").Value = 38.76511933333333
$ws.Range("N4").Value = 116.295358
$ws.Range("O4").Value = 0.3601657659886822
$ws.Range("P4").Value = 0.3601657659886822
$ws.Range("Q4").Value = 19.65454866561578
$ws.Range("R4").Value = 176.890937990542
$ws.Range("S4").Value = 0.2511692527370455
$ws.Range("T4").Value = 0.2511692527370455

# Row 5
$ws.Range("H5").Value = 0.524222
$ws.Range("I5").Value = 0.2403455475586373
$ws.Range("J5").Value = 0.2403455475586373
$ws.Range("M5").Value = 0.07904133333333334
$ws.Range("N5").Value = 0.237124
$ws.Range("O5").Value = 0.0007343710751920149
$ws.Range("P5").Value = 0.0007343710751920148
$ws.Range("Q5").Value = 0.01381173528088889
$ws.Range("R5").Value = 0.124305617528
$ws.Range("S5").Value = 0.00017650281817825
$ws.Range("T5").Value = 0.00017650281817825

# Row 6
$ws.Range("H6").Value = 0.524222
$ws.Range("I6").Value = 0.2403455475586373
$ws.Range("J6").Value = 0.2403455475586373
$ws.Range("M6").Value = 68.78716633333333
$ws.Range("N6").Value = 206.361499
$ws.Range("O6").Value = 0.6390998629361258
$ws.Range("P6").Value = 0.6390998629361259
$ws.Range("Q6").Value = 12.01991530319756
$ws.Range("R6").Value = 108.179237728778
$ws.Range("S6").Value = 0.1536048065020332
$ws.Range("T6").Value = 0.1536048065020332

# Row 7
$ws.Range("H7").Value = 0.524222
$ws.Range("I7").Value = 0.2403455475586373
$ws.Range("J7").Value = 0.2403455475586373
$ws.Range("M7").Value = 38.76511933333333
$ws.Range("N7").Value = 116.295358
$ws.Range("O7").Value = 0.3601657659886822
$ws.Range("P7").Value = 0.3601657659886822
$ws.Range("Q7").Value = 6.773842795719554
$ws.Range("R7").Value = 60.96458516147599
$ws.Range("S7").Value = 0.08656423823842586
$ws.Range("T7").Value = 0.08656423823842586

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("G8").Value = 0.04528233333333333
$ws.Range("H8").Value = 0.135847
$ws.Range("I8").Value = 0.06228319604899872
$ws.Range("J8").Value = 0.06228319604899872
$ws.Range("M8").Value = 0.07904133333333334
$ws.Range("N8").Value = 0.237124
$ws.Range("O8").Value = 0.0007343710751920149
$ws.Range("P8").Value = 0.0007343710751920148
$ws.Range("Q8").Value = 0.003579176003111112
$ws.Range("R8").Value = 0.032212584028
$ws.Range("S8").Value = 0.00004573897764889825
$ws.Range("T8").Value = 0.00004573897764889824

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("G9").Value = 0.04528233333333333
$ws.Range("H9").Value = 0.135847
$ws.Range("I9").Value = 0.06228319604899872
$ws.Range("J9").Value = 0.06228319604899872
$ws.Range("M9").Value = 68.78716633333333
$ws.Range("N9").Value = 206.361499
$ws.Range("O9").Value = 0.6390998629361258
$ws.Range("P9").Value = 0.6390998629361259
$ws.Range("Q9").Value = 3.114843394961444
$ws.Range("R9").Value = 28.033590554653
$ws.Range("S9").Value = 0.03980518205813893
$ws.Range("T9").Value = 0.03980518205813894

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("G10").Value = 0.04528233333333333
$ws.Range("H10").Value = 0.135847
$ws.Range("I10").Value = 0.06228319604899872
$ws.Range("J10").Value = 0.06228319604899872
$ws.Range("M10").Value = 38.76511933333333
$ws.Range("N10").Value = 116.295358
$ws.Range("O10").Value = 0.3601657659886822
$ws.Range("P10").Value = 0.3601657659886822
$ws.Range("Q10").Value = 1.755375055358444
$ws.Range("R10").Value = 15.798375498226
$ws.Range("S10").Value = 0.02243227501321089
$ws.Range("T10").Value = 0.02243227501321089

